$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently runs through column J (year 2022). Extend it one
# column to the right with a 2023 column, copying the formatting that
# column J already has (header style + the three data-row styles).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the 2023 figures.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 522.6
$ws.Range("K5").Value = 377.8
$ws.Range("K6").Value = 661
